$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp refresh ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 12:52"

# --- Data refresh + re-sort for affected country rows ---

# Espana (row 5)
$ws.Range("B5").Value = 232128
$ws.Range("C5").Value = 2706
$ws.Range("D5").Value = 123903
$ws.Range("E5").Value = 84403
$ws.Range("G5").Value = 301
$ws.Range("H5").Value = 23822

# Suiza (row 19)
$ws.Range("B19").Value = 29264
$ws.Range("C19").Value = 100
$ws.Range("E19").Value = 5387

# Uzbekistan (row 68)
$ws.Range("D68").Value = 934
$ws.Range("E68").Value = 997

# Malta (row 111)
$ws.Range("B111").Value = 458
$ws.Range("C111").Value = 8
$ws.Range("D111").Value = 303
$ws.Range("E111").Value = 151

# Rows 129/130 swap rank: Maldivas overtakes Paraguay with updated totals
$ws.Range("A129").Value = "Maldivas"
$ws.Range("B129").Value = 245
$ws.Range("C129").Value = 19
$ws.Range("D129").Value = 17
$ws.Range("E129").Value = 228
$ws.Range("F129").Value = 2
$ws.Range("H129").Value = 0

$ws.Range("A130").Value = "Paraguay"
$ws.Range("B130").Value = 230
$ws.Range("C130").Value = 2
$ws.Range("D130").Value = 95
$ws.Range("E130").Value = 126
$ws.Range("F130").Value = 1
$ws.Range("H130").Value = 9

# Rows 151/152 swap rank: Zambia overtakes Monaco with updated totals
$ws.Range("A151").Value = "Zambia"
$ws.Range("C151").Value = 6
$ws.Range("E151").Value = 50
$ws.Range("H151").Value = 3

$ws.Range("A152").Value = "Monaco"
$ws.Range("B152").Value = 95
$ws.Range("E152").Value = 49
$ws.Range("H152").Value = 4

# Nepal (row 167)
$ws.Range("B167").Value = 54
$ws.Range("C167").Value = 2
$ws.Range("E167").Value = 38
